# #220 adapt circulating supply
# Append the new "2.5.1" changelog entry as row 39 of the (only) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New changelog row: Date | Version | Changes
$ws.Range("A39").Value = 44648
$ws.Range("B39").Value = "2.5.1"
$ws.Range("C39").Value = "Improvement:`n- Adapt circulating supply in Overview and Coins evaluation"

# Match formatting of the existing changelog rows: wrapped, two-line text at
# row height 30 (same as every other 2-line entry, e.g. row 2/3/10/11/21).
$ws.Range("C39").WrapText = $true
$ws.Rows.Item(39).RowHeight = 30

# The workbook was left with this cell selected/active.
$null = $ws.Range("J38").Select()
